# Generate Report for Handback
# - Removes the "d0c6ac18-ceaf-4b0d-9442-1483750bf9e2" row (row 3) from every
#   sheet (Overview, zh-cn, de-de), including its now-stale hyperlinks.
# - Refreshes the "Correspond Handoff/Handback" timestamps on the remaining
#   (7cfe3a7e-...) row for the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

function Remove-RowAndHyperlinks($ws, $row) {
    # Drop any hyperlinks anchored on the row we are about to remove -
    # EntireRow delete leaves dangling <hyperlink> entries behind otherwise.
    $changed = $true
    while ($changed) {
        $changed = $false
        foreach ($hl in $ws.Hyperlinks) {
            if ($hl.Range.Row -eq $row) {
                $hl.Delete()
                $changed = $true
                break
            }
        }
    }
    $ws.Rows.Item($row).Delete()
}

# --- Overview sheet: drop row 3 (d0c6ac18-...) ---
$wsOverview = $wb.Worksheets.Item("Overview")
Remove-RowAndHyperlinks $wsOverview 3

# --- zh-cn sheet: refresh row 2 timestamps, drop row 3 ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Cells.Item(2, 5).Value = "2016-03-11 12:31:48"
$wsZh.Cells.Item(2, 8).Value = "2016-03-11 12:32:05"
Remove-RowAndHyperlinks $wsZh 3

# --- de-de sheet: refresh row 2 timestamps, drop row 3 ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Cells.Item(2, 5).Value = "2016-03-11 12:31:51"
$wsDe.Cells.Item(2, 8).Value = "2016-03-11 12:32:11"
Remove-RowAndHyperlinks $wsDe 3
